# Add three new match columns (HD, HE, HF) to the Essendon stats sheet,
# matching the formatting pattern already used by the table: every data
# column carries style index 1 except the very last (newest) column, which
# stays unstyled until a further column is appended after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: (row, HC value [unchanged], HD value, HE value, HF value)
$rows = @(
    @(1, 10239, 10251, 10259, 10266),
    @(2, 2020, 2020, 2020, 2020),
    @(3, 9, 10, 11, 12),
    @(4, 0, 0, 1, 1),
    @(5, 1, 1, 0.5, 0),
    @(6, 28, 55, 73, 33),
    @(7, 91, 59, 73, 68),
    @(8, -63, -4, 0, -35),
    @(9, 0, 0, 0.5, 0),
    @(10, 2, 9, 8, 15),
    @(11, 140, 183, 177, 172),
    @(12, 108, 124, 140, 135),
    @(13, 248, 307, 317, 307),
    @(14, 1.3, 1.48, 1.26, 1.27),
    @(15, 41, 69, 74, 80),
    @(16, 39, 44, 44, 37),
    @(17, 34, 22, 24, 22),
    @(18, 16, 19, 16, 11),
    @(19, 24, 26, 12, 11),
    @(20, 3, 8, 11, 5),
    @(21, 0, 5, 9, 5),
    @(22, 9, 5, 5, 3),
    @(23, 1, 2, 2, 0),
    @(24, 13, 15, 18, 8),
    @(25, 23.1, 53.3, 61.1, 62.5),
    @(26, 82.67, 38.38, 28.82, 61.4),
    @(27, 19.08, 20.47, 17.61, 38.38),
    @(28, 34, 36, 26, 30),
    @(29, 46, 64, 46, 47),
    @(30, 23, 39, 30, 27),
    @(31, 42, 39, 48, 42),
    @(32, 3.23, 2.6, 2.67, 5.25),
    @(33, 14, 4.88, 4.36, 8.4),
    @(34, 28.6, 33.3, 33.3, 19),
    @(35, 7.1, 20.5, 22.9, 11.9),
    @(36, 186, 186.1, 186, 187),
    @(37, 85.5, 86.3, 85.9, 87.5),
    @(38, 25, 25.16, 25.16, 24.66),
    @(39, 77.1, 79.7, 79.9, 67.2),
    @(40, 9, 8, 8, 9),
    @(41, 6, 8, 8, 8),
    @(42, 4, 3, 4, 3),
    @(43, 3, 3, 2, 2),
    @(44, 110, 120, 117, 90),
    @(45, 133, 173, 188, 209),
    @(46, 169, 200, 230, 236),
    @(47, 68.1, 65.1, 72.6, 76.9),
    @(48, 46, 64, 46, 47),
    @(49, 5, 4, 12, 3),
    @(50, 5, 8, 13, 9),
    @(51, 34, 36, 26, 30),
    @(52, 23, 39, 30, 27),
    @(53, 35, 44, 34, 39),
    @(54, 6, 5, 7, 7),
    @(55, 0, 5, 9, 5),
    @(56, 0, 62.5, 81.8, 100),
    @(57, 210, 175, 165, 191),
    @(58, 107, 131, 82, 107),
    @(59, 317, 306, 247, 298),
    @(60, 1.96, 1.34, 2.01, 1.79),
    @(61, 99, 55, 64, 98),
    @(62, 50, 55, 42, 38),
    @(63, 25, 39, 32, 34),
    @(64, 24, 26, 12, 11),
    @(65, 16, 19, 16, 11),
    @(66, 14, 8, 11, 10),
    @(67, 11, 6, 6, 6),
    @(68, 6, 7, 6, 6),
    @(69, 1, 4, 1, 2),
    @(70, 21, 19, 18, 18),
    @(71, 66.7, 42.1, 61.1, 55.6),
    @(72, 22.64, 38.25, 22.45, 29.8),
    @(73, 15.1, 16.11, 13.72, 16.56),
    @(74, 29, 36, 32, 25),
    @(75, 36, 51, 40, 42),
    @(76, 38, 31, 36, 37),
    @(77, 37, 47, 43, 37),
    @(78, 1.76, 2.47, 2.39, 2.06),
    @(79, 2.64, 5.88, 3.91, 3.7),
    @(80, 54.1, 31.9, 39.5, 43.2),
    @(81, 37.8, 17, 25.6, 27),
    @(82, 188.6, 188.2, 187.3, 187.8),
    @(83, 88.6, 87.7, 85.5, 84.2),
    @(84, 24.8, 25.74, 23.8, 24.66),
    @(85, 85.1, 104.3, 70.5, 77.9),
    @(86, 7, 7, 10, 10),
    @(87, 7, 5, 6, 5),
    @(88, 4, 5, 4, 4),
    @(89, 4, 5, 2, 3),
    @(90, 125, 114, 114, 107),
    @(91, 195, 184, 126, 192),
    @(92, 245, 204, 160, 208),
    @(93, 77.3, 66.7, 64.8, 69.8),
    @(94, 36, 51, 40, 42),
    @(95, 9, 4, 7, 8),
    @(96, 15, 7, 10, 18),
    @(97, 29, 36, 32, 25),
    @(98, 38, 31, 36, 37),
    @(99, 32, 59, 39, 41),
    @(100, 2, 1, 2, 11),
    @(101, 11, 6, 6, 6),
    @(102, 78.6, 75, 54.5, 60)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $hcValue = $r[1]
    $hdValue = $r[2]
    $heValue = $r[3]
    $hfValue = $r[4]

    # HF inherits HC's current (unstyled) formatting before HC itself is restyled,
    # so the newly-appended last column ends up unstyled, like HC used to be.
    $ws.Range("HC" + $rowNum).Copy()
    $ws.Range("HF" + $rowNum).PasteSpecial(-4122)

    # New cells pick up style 1 automatically from their styled left neighbor.
    $ws.Range("HD" + $rowNum).Value = $hdValue
    $ws.Range("HE" + $rowNum).Value = $heValue
    $ws.Range("HF" + $rowNum).Value = $hfValue

    # Re-enter HC as if it were freshly typed so it also inherits style 1
    # from its (now styled) neighbors, matching the rest of the row.
    $ws.Range("HC" + $rowNum).ClearContents()
    $ws.Range("HC" + $rowNum).Value = $hcValue
}

$excel.CutCopyMode = 0

Write-Host "Updated columns HD:HF for rows 1-102; new dimension should be A1:HF102"
